$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.98071881606765199
$ws.Range("I2").NumberFormat = $ws.Range("I3").NumberFormat

$ws.Range("I3").Value = 0.98071881606765199
$ws.Range("I4").Value = 0.99076109936574897
$ws.Range("I5").Value = 0.99298097251585504
$ws.Range("I6").Value = 0.947779021426287
$ws.Range("I7").Value = 0.97173968660057497
$ws.Range("I8").Value = 0.97804285257435197
$ws.Range("I9").Value = 0.98357211384713805
$ws.Range("I10").Value = 0.72402640264026397
$ws.Range("I11").Value = 0.71724186704384696
$ws.Range("I12").Value = 0.829434229137199
$ws.Range("I13").Value = 0.86489391796322501
$ws.Range("I14").Value = 0.71525763754499905
$ws.Range("I15").Value = 0.78692712279826305
$ws.Range("I16").Value = 0.85457085289326196
$ws.Range("I17").Value = 0.87163815302556102

$ws.Range("J7").Select()
